$wb = $excel.ActiveWorkbook

# The same three rows are updated identically on both the "展览" sheet
# and the "全部类型" sheet (which aggregates all event types).
$sheetNames = @("展览", "全部类型")

foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)

    # F2: 想去人数 (number of people interested) 568 -> 570
    $ws.Range("F2").Value = 570

    # F9: 想去人数 3598 -> 3606
    $ws.Range("F9").Value = 3606

    # F10: 想去人数 59 -> 60
    $ws.Range("F10").Value = 60
}
